$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update average_county_temperature (column AA) values for the facilities
# that received updated NOAA temperature data. Values are grouped below by
# the contiguous worksheet rows that share the same facility/value.
$ws.Range("AA2:AA10").Value = 19.30324074074072
$ws.Range("AA11:AA16").Value = 17.25771604938272
$ws.Range("AA17:AA25").Value = 13.62268518518517
$ws.Range("AA65:AA67").Value = 13.75752314814816
$ws.Range("AA68:AA73").Value = 19.79629629629628
$ws.Range("AA74:AA79").Value = 0.8611111111111096
$ws.Range("AA83:AA88").Value = 0.8611111111111096
$ws.Range("AA92:AA94").Value = 5.486111111111112
$ws.Range("AA95:AA100").Value = 16.86342592592595
$ws.Range("AA101:AA106").Value = 5.486111111111112
$ws.Range("AA107:AA109").Value = 12.41429539295394
$ws.Range("AA113:AA115").Value = 19.36574074074073
$ws.Range("AA140:AA142").Value = 12.41429539295394
$ws.Range("AA146:AA154").Value = 12.41429539295394
